$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69; existing rows 69..138 shift down to 70..139
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new record
$ws.Cells.Item(69, 1).Value = 1
$ws.Cells.Item(69, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(69, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(69, 4).Value = 44874
$ws.Cells.Item(69, 5).Value = 15
$ws.Cells.Item(69, 6).Value = 'Fruta'
$ws.Cells.Item(69, 7).Value = 100106
$ws.Cells.Item(69, 8).Value = 'Oleaginosos'
$ws.Cells.Item(69, 9).Value = 100106002
$ws.Cells.Item(69, 10).Value = 'Palta'
$ws.Cells.Item(69, 11).Value = 'Fuerte'
$ws.Cells.Item(69, 12).Value = 'Tercera'
$ws.Cells.Item(69, 13).Value = 200
$ws.Cells.Item(69, 14).Value = 53000
$ws.Cells.Item(69, 15).Value = 55000
$ws.Cells.Item(69, 16).Value = 54000
$ws.Cells.Item(69, 17).Value = '$/caja 25 kilos'
$ws.Cells.Item(69, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(69, 19).Value = 2160
$ws.Cells.Item(69, 20).Value = 25
